$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shotgun's "fire delay" parameter was moved into the playerclass,
# so this row's key is renamed from SHOTGUN_FIRE_DELAY to SHOTGUN_BULLET_RELOAD.
$ws.Range("A7").Value = "SHOTGUN_BULLET_RELOAD"

# Leave selection on the edited cell, matching the saved view state.
$ws.Range("A7").Select()
